# "Bug fixes to export_metadata"
#
# The reviewer's H40:H51 scratch notes (a "DELETE" / "sample-type_group"
# checklist) are cut out of Sheet1 and moved onto a brand-new sheet named
# "instructions for Conrad" placed right after Sheet1. The three cells
# that held the stray "All gear" placeholder text are cleared out, which
# lets the now-unreferenced shared string get garbage-collected on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. Create the destination sheet right after Sheet1.
$notes = $wb.Worksheets.Add($null, $ws)
$notes.Name = "instructions for Conrad"

# 2. Move (copy values + formats) the H40:H51 note block onto the new
#    sheet's A1:A12 before anything on Sheet1 is touched.
$ws.Range("H40:H51").Copy()
$notes.Range("A1:A12").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("H40:H51").Copy()
$notes.Range("A1:A12").PasteSpecial(-4122)   # xlPasteFormats
$notes.Range("C38").Select()

# 3. Clear the three stray "All gear" cells on Sheet1.
$ws.Range("D16").ClearContents()
$ws.Range("D21").ClearContents()
$ws.Range("D35").ClearContents()

# 4. Remove the now-empty H40:H51 rows from Sheet1 entirely.
$ws.Range("A40:H51").EntireRow.Delete()

# 5. Restore Sheet1 as the active/selected sheet and selection.
$ws.Activate()
$ws.Range("D24").Select()
